{"js": "// Fix typo in example code: `parameters[` -> `parameter[`\n// (the function signature uses the singular `parameter`, the body had a\n// leftover plural typo in four places across two code snippets).\nconst results = context.document.body.search(\"parameters[\", {\n  matchCase: true,\n  matchWholeWord: false,\n});\nresults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  const item = results.items[i];\n  // item.text is the literal matched text \"parameters[\" \u2014 replace just\n  // that occurrence in place, preserving the rest of the run/paragraph.\n  item.insertText(\"parameter[\", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Fix typo in example: `parameters[` -> `parameter[`\n# The function signature uses the singular \"parameter\"; the body had a\n# leftover plural typo (\"parameters[\") in four spots across two code\n# snippets (capital_cost and fixed_cost).\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"parameters[\"\n$find.Replacement.Text = \"parameter[\"\n\n$find.Execute(\n    $find.Text,               # FindText\n    $true,                    # MatchCase\n    $false,                   # MatchWholeWord\n    $false,                   # MatchWildcards\n    $false,                   # MatchSoundsLike\n    $false,                   # MatchAllWordForms\n    $true,                    # Forward\n    1,                        # Wrap: wdFindContinue\n    $false,                   # Format\n    $find.Replacement.Text,   # ReplaceWith\n    2                         # Replace: wdReplaceAll\n) | Out-Null\n"}
